$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text while we write the new
# values, so strings that happen to look like numbers (e.g. "1.004") are not
# silently converted to numeric cells. Restore formatting afterwards so no
# stray number-format is left applied to any cell.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '30.301.75'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '1.928.16'
$ws.Range('E3').Value = '  -0.09%  '

$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.34%  '

$ws.Range('D5').Value = '0.7396'
$ws.Range('E5').Value = '  +2.72%  '

$ws.Range('D6').Value = '243.45'
$ws.Range('E6').Value = '  -2.48%  '

$ws.Range('E7').Value = '  +0.29%  '

$ws.Range('D8').Value = '27.46'
$ws.Range('E8').Value = '  -1.67%  '

$ws.Range('D9').Value = '0.3147'
$ws.Range('E9').Value = '  -1.74%  '

$ws.Range('D10').Value = '0.07017'
$ws.Range('E10').Value = '  -1.15%  '

$ws.Range('D11').Value = '0.08033'
$ws.Range('E11').Value = '  +0.23%  '

$ws.Range('D12').Value = '0.7743'
$ws.Range('E12').Value = '  -1.76%  '

$ws.Range('D13').Value = '1.911.61'
$ws.Range('E13').Value = '  -0.97%  '

$ws.Range('D14').Value = '5.361'
$ws.Range('E14').Value = '  -0.30%  '

$ws.Range('D15').Value = '92.97'

$ws.Range('D16').Value = '14.43'
$ws.Range('E16').Value = '  -1.50%  '

$ws.Range('D17').Value = '30.318.83'
$ws.Range('E17').Value = '  +0.08%  '

$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '249.85'
$ws.Range('E18').Value = '  -2.71%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '5.943'
$ws.Range('E19').Value = '  +3.73%  '

$ws.Range('D20').Value = '0.000007927'
$ws.Range('E20').Value = '  -1.62%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.192.38'
$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.29%  '

$ws.Range('D23').Value = '1.004'
$ws.Range('E23').Value = '  +0.55%  '

$ws.Range('D24').Value = '6.635'
$ws.Range('E24').Value = '  -2.66%  '

$ws.Range('D25').Value = '9.540'
$ws.Range('E25').Value = '  -0.13%  '

$ws.Range('D26').Value = '165.62'
$ws.Range('E26').Value = '  +0.69%  '

$ws.Range('D27').Value = '18.99'
$ws.Range('E27').Value = '  -0.54%  '

$ws.Range('D28').Value = '0.1281'
$ws.Range('E28').Value = '  -0.03%  '

$ws.Range('D29').Value = '2.161'
$ws.Range('E29').Value = '  -5.89%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.357'
$ws.Range('E30').Value = '  +0.31%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.559'
$ws.Range('E31').Value = '  +1.67%  '

$ws.Range('D32').Value = '4.401'
$ws.Range('E32').Value = '  -0.43%  '

$ws.Range('D33').Value = '4.102'
$ws.Range('E33').Value = '  -1.19%  '

$ws.Range('D34').Value = '0.05208'
$ws.Range('E34').Value = '  +1.90%  '

$ws.Range('D35').Value = '1.306'
$ws.Range('E35').Value = '  +1.37%  '

$ws.Range('D36').Value = '0.7533'
$ws.Range('E36').Value = '  +0.51%  '

$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  +0.03%  '

$ws.Range('D38').Value = '0.01946'
$ws.Range('E38').Value = '  -1.98%  '

$ws.Range('D39').Value = '2.790'
$ws.Range('E39').Value = '  -0.31%  '

$ws.Range('D40').Value = '6.516'
$ws.Range('E40').Value = '  +1.78%  '

$ws.Range('D41').Value = '76.47'
$ws.Range('E41').Value = '  -2.48%  '

$ws.Range('D42').Value = '0.4491'
$ws.Range('E42').Value = '  -0.69%  '

$ws.Range('D43').Value = '1.948'
$ws.Range('E43').Value = '  -2.39%  '

$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.27%  '

$ws.Range('D45').Value = '0.8388'
$ws.Range('E45').Value = '  -0.81%  '

$ws.Range('D46').Value = '7.669'
$ws.Range('E46').Value = '  +2.45%  '

$ws.Range('D47').Value = '9.934'
$ws.Range('E47').Value = '  +0.94%  '

$ws.Range('D48').Value = '101.22'
$ws.Range('E48').Value = '  -0.12%  '

$ws.Range('D49').Value = '37.43'
$ws.Range('E49').Value = '  +1.57%  '

$ws.Range('D50').Value = '2.063.79'
$ws.Range('E50').Value = '  -1.39%  '

$ws.Range('D51').Value = '0.1225'
$ws.Range('E51').Value = '  +7.13%  '

# Restore the Price column formatting to General (matching the workbook
# defaults) now that the text has been written.
$priceRange.ClearFormats()